$wb = $excel.ActiveWorkbook

# Rename the "Dragons" sheet to "味全龍"
$dragonsSheet = $wb.Worksheets.Item("Dragons")
$dragonsSheet.Name = "味全龍"

# Make it the active/selected sheet (was previously "台鋼雄鷹")
$dragonsSheet.Activate()
$dragonsSheet.Select()
